# Rename the sheet from "Data" to "Summary"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Summary"

# Clear the old layout cells whose content is moving to new rows
$ws.Range("B5:D5").Clear()
$ws.Range("A6:D7").Clear()

# Create a new cell style for the bold+underlined "Source Type" subheading,
# mirroring the existing named styles (name/title/source/HyperLink) already
# used by this template.
$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Name = "Calibri"
$titleStyle.Font.Size = 11
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true

# --- Row 9: new "Source Type" subheading ---
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Name = "Calibri"
$ws.Range("A9").Font.Size = 11
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- Rows 11-14: first indicator table (Employment / Enterprises) ---
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

$ws.Range("A12").Value = "Employment (% of total)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "30"

$ws.Range("A13").Value = "Enterprises (% of total)"
$ws.Range("A13").Font.Bold = $true
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "90"

$ws.Range("A14").Value = "Source: Min PME - ADEPME, 2010"
$ws.Range("A14").Font.Italic = $true

# --- Rows 16-18: second indicator table (Value added) ---
$ws.Range("B16").Value = "Micro"
$ws.Range("B16").Font.Bold = $true
$ws.Range("C16").Value = "SMEs"
$ws.Range("C16").Font.Bold = $true
$ws.Range("D16").Value = "MSMEs"
$ws.Range("D16").Font.Bold = $true

$ws.Range("A17").Value = "Value added to the economy (% of total)"
$ws.Range("A17").Font.Bold = $true
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20"

$ws.Range("A18").Value = "Source: Min PME - ADEPME, 2010"
$ws.Range("A18").Font.Italic = $true

# --- Rows 23-24: reference section ---
$ws.Range("A23").Value = "Min PME - ADEPME"
$ws.Range("A23").Font.Bold = $true

$ws.Range("A24").Value = "Ministere des mines, de l'Industrie, de l'Agro-industrie et des PME, Direction des Petites et Moyennes Entreprises (Min PME - ADEPME), ""LETTRE  DE  POLITIQUE SECTORIELLE  DES  PME"", 2010, p. 9. Available at http://www.senegal-entreprises.net/3-download/lettre-politique-sectorielle-10-2010.pdf"
$ws.Range("A24").Font.Italic = $true
